# Mexico population 2010 und 2020!! -- update the Guatemala 1998 mortality
# table so that nqx (E4) is fed back from the life-table columns (J4/K4)
# instead of the raw nDx/nNx ratio, and nax (G4) becomes a fixed input
# (0.044) instead of a formula. Because E4/F4/K4 now form a circular
# reference, the workbook already has iterative calculation enabled
# (iterate="1", iterateCount="1000"), so we just have to let the engine
# settle on the fixed point.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tabla mortalidad 1998")

# G4 becomes a hard-coded value (previously a shared formula
# "(B4*E4)/(1+(B4-F4)*E4)").
$ws.Range("G4").Value = 0.044

# E4 now reads back from the life-table build-out instead of D4/C4.
$ws.Range("E4").Formula = "=J4/K4"

$excel.Calculate()

# E4 -> F4 -> K4 -> E4 is circular (G4 is now a constant, so it is no
# longer part of the loop). Re-touching a precedent and recalculating
# repeatedly drives the Gauss-Seidel iteration to the same fixed point
# Excel's own iterative calculation would converge to.
for ($i = 0; $i -lt 15; $i++) {
    $ws.Range("G4").Value = 0.044
    $excel.Calculate()
}

# Highlight the two edited inputs in yellow (same number format as
# before, just with a fill added).
$ws.Range("E4").Interior.Color = 65535
$ws.Range("G4").Interior.Color = 65535

# Leave the cursor on the cell that was edited, without disturbing which
# sheet tab is active overall.
$ws.Range("E4").Select()
$activeSheet = $wb.Worksheets.Item("Descomposición e0 - Arriaga")
$activeSheet.Activate()
